$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump "Forandrad" (column C) date value for every data row 2-37
$ws.Range("C2:C37").Value = 46081

# Row 5  <-  old row 6 data (A 21219-2023)
$ws.Range("A5").Value = "A 21219-2023"
$ws.Range("B5").Value = 45062.0
$ws.Range("G5").Value = 2.6
$ws.Range("H5").Value = 7.0
$ws.Range("I5").Value = 11.0
$ws.Range("J5").Value = 1.0
$ws.Range("K5").Value = 0.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 0.0
$ws.Range("N5").Value = 0.0
$ws.Range("O5").Value = 2.0
$ws.Range("P5").Value = 1.0
$ws.Range("Q5").Value = 16.0
$ws.Range("R5").Value = "Ask`r`nOlivbrun spindling`r`nGrov baronmossa`r`nMurgröna`r`nNästrot`r`nPlatt fjädermossa`r`nPurpurknipprot`r`nSkogsknipprot`r`nStrimspindling`r`nSårläka`r`nTraslav`r`nTvåblad`r`nUnderviol`r`nFläcknycklar`r`nBlåsippa`r`nGullviva"
$ws.Range("S5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/artfynd/A 21219-2023 artfynd.xlsx`", `"A 21219-2023`")"
$ws.Range("T5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/kartor/A 21219-2023 karta.png`", `"A 21219-2023`")"
$ws.Range("V5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomål/A 21219-2023 FSC-klagomål.docx`", `"A 21219-2023`")"
$ws.Range("W5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomålsmail/A 21219-2023 FSC-klagomål mail.docx`", `"A 21219-2023`")"
$ws.Range("X5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsyn/A 21219-2023 tillsynsbegäran.docx`", `"A 21219-2023`")"
$ws.Range("Y5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsynsmail/A 21219-2023 tillsynsbegäran mail.docx`", `"A 21219-2023`")"

# Row 6  <-  old row 5 data (A 27865-2024)
$ws.Range("A6").Value = "A 27865-2024"
$ws.Range("B6").Value = 45475.67303240741
$ws.Range("G6").Value = 0.7
$ws.Range("H6").Value = 5.0
$ws.Range("I6").Value = 8.0
$ws.Range("J6").Value = 6.0
$ws.Range("K6").Value = 0.0
$ws.Range("L6").Value = 0.0
$ws.Range("M6").Value = 0.0
$ws.Range("N6").Value = 0.0
$ws.Range("O6").Value = 6.0
$ws.Range("P6").Value = 0.0
$ws.Range("Q6").Value = 16.0
$ws.Range("R6").Value = "Bergjohannesört`r`nBredarun`r`nKlasefibbla`r`nSolvända`r`nSvart taggsvamp`r`nVippärt`r`nBlå slemspindling`r`nDiskvaxskivling`r`nFjällig taggsvamp s.str.`r`nMurgröna`r`nNästrot`r`nSkogsknipprot`r`nSträvlosta`r`nUnderviol`r`nGrönvit nattviol`r`nBlåsippa"
$ws.Range("S6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/artfynd/A 27865-2024 artfynd.xlsx`", `"A 27865-2024`")"
$ws.Range("T6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/kartor/A 27865-2024 karta.png`", `"A 27865-2024`")"
$ws.Range("V6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomål/A 27865-2024 FSC-klagomål.docx`", `"A 27865-2024`")"
$ws.Range("W6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomålsmail/A 27865-2024 FSC-klagomål mail.docx`", `"A 27865-2024`")"
$ws.Range("X6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsyn/A 27865-2024 tillsynsbegäran.docx`", `"A 27865-2024`")"
$ws.Range("Y6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsynsmail/A 27865-2024 tillsynsbegäran mail.docx`", `"A 27865-2024`")"

# Row 10  <-  old row 12 data (A 24778-2022)
$ws.Range("A10").Value = "A 24778-2022"
$ws.Range("B10").Value = 44728.0
$ws.Range("G10").Value = 1.0
$ws.Range("H10").Value = 0.0
$ws.Range("I10").Value = 2.0
$ws.Range("J10").Value = 3.0
$ws.Range("K10").Value = 0.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 0.0
$ws.Range("N10").Value = 0.0
$ws.Range("O10").Value = 4.0
$ws.Range("P10").Value = 1.0
$ws.Range("Q10").Value = 6.0
$ws.Range("R10").Value = "Tofsäxing`r`nFlentimotej`r`nPoppeltofsskivling`r`nSlåtterfibbla`r`nKalktallört`r`nMurgröna"
$ws.Range("S10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/artfynd/A 24778-2022 artfynd.xlsx`", `"A 24778-2022`")"
$ws.Range("T10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/kartor/A 24778-2022 karta.png`", `"A 24778-2022`")"
$ws.Range("V10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomål/A 24778-2022 FSC-klagomål.docx`", `"A 24778-2022`")"
$ws.Range("W10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomålsmail/A 24778-2022 FSC-klagomål mail.docx`", `"A 24778-2022`")"
$ws.Range("X10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsyn/A 24778-2022 tillsynsbegäran.docx`", `"A 24778-2022`")"
$ws.Range("Y10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsynsmail/A 24778-2022 tillsynsbegäran mail.docx`", `"A 24778-2022`")"

# Row 11  <-  old row 10 data (A 27636-2023)
$ws.Range("A11").Value = "A 27636-2023"
$ws.Range("B11").Value = 45097.0
$ws.Range("G11").Value = 7.4
$ws.Range("H11").Value = 2.0
$ws.Range("I11").Value = 1.0
$ws.Range("J11").Value = 2.0
$ws.Range("K11").Value = 0.0
$ws.Range("L11").Value = 1.0
$ws.Range("M11").Value = 0.0
$ws.Range("N11").Value = 0.0
$ws.Range("O11").Value = 3.0
$ws.Range("P11").Value = 1.0
$ws.Range("Q11").Value = 6.0
$ws.Range("R11").Value = "Ask`r`nBacktimjan`r`nSolvända`r`nMurgröna`r`nAlvarmalört`r`nGullviva"
$ws.Range("S11").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/artfynd/A 27636-2023 artfynd.xlsx`", `"A 27636-2023`")"
$ws.Range("T11").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/kartor/A 27636-2023 karta.png`", `"A 27636-2023`")"
$ws.Range("V11").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomål/A 27636-2023 FSC-klagomål.docx`", `"A 27636-2023`")"
$ws.Range("W11").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomålsmail/A 27636-2023 FSC-klagomål mail.docx`", `"A 27636-2023`")"
$ws.Range("X11").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsyn/A 27636-2023 tillsynsbegäran.docx`", `"A 27636-2023`")"
$ws.Range("Y11").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsynsmail/A 27636-2023 tillsynsbegäran mail.docx`", `"A 27636-2023`")"

# Row 12  <-  old row 11 data (A 38039-2022)
$ws.Range("A12").Value = "A 38039-2022"
$ws.Range("B12").Value = 44811.0
$ws.Range("G12").Value = 7.4
$ws.Range("H12").Value = 2.0
$ws.Range("I12").Value = 1.0
$ws.Range("J12").Value = 2.0
$ws.Range("K12").Value = 0.0
$ws.Range("L12").Value = 1.0
$ws.Range("M12").Value = 0.0
$ws.Range("N12").Value = 0.0
$ws.Range("O12").Value = 3.0
$ws.Range("P12").Value = 1.0
$ws.Range("Q12").Value = 6.0
$ws.Range("R12").Value = "Ask`r`nBacktimjan`r`nSolvända`r`nMurgröna`r`nAlvarmalört`r`nGullviva"
$ws.Range("S12").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/artfynd/A 38039-2022 artfynd.xlsx`", `"A 38039-2022`")"
$ws.Range("T12").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/kartor/A 38039-2022 karta.png`", `"A 38039-2022`")"
$ws.Range("V12").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomål/A 38039-2022 FSC-klagomål.docx`", `"A 38039-2022`")"
$ws.Range("W12").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomålsmail/A 38039-2022 FSC-klagomål mail.docx`", `"A 38039-2022`")"
$ws.Range("X12").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsyn/A 38039-2022 tillsynsbegäran.docx`", `"A 38039-2022`")"
$ws.Range("Y12").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsynsmail/A 38039-2022 tillsynsbegäran mail.docx`", `"A 38039-2022`")"

# Row 18  <-  old row 19 data (A 2864-2026)
$ws.Range("A18").Value = "A 2864-2026"
$ws.Range("B18").Value = 46038.0
$ws.Range("G18").Value = 1.0
$ws.Range("H18").Value = 2.0
$ws.Range("I18").Value = 1.0
$ws.Range("J18").Value = 1.0
$ws.Range("K18").Value = 0.0
$ws.Range("L18").Value = 0.0
$ws.Range("M18").Value = 0.0
$ws.Range("N18").Value = 0.0
$ws.Range("O18").Value = 1.0
$ws.Range("P18").Value = 0.0
$ws.Range("Q18").Value = 3.0
$ws.Range("R18").Value = "Spillkråka`r`nMindre märgborre`r`nBlåsippa"
$ws.Range("S18").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/artfynd/A 2864-2026 artfynd.xlsx`", `"A 2864-2026`")"
$ws.Range("T18").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/kartor/A 2864-2026 karta.png`", `"A 2864-2026`")"
$ws.Range("V18").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomål/A 2864-2026 FSC-klagomål.docx`", `"A 2864-2026`")"
$ws.Range("W18").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomålsmail/A 2864-2026 FSC-klagomål mail.docx`", `"A 2864-2026`")"
$ws.Range("X18").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsyn/A 2864-2026 tillsynsbegäran.docx`", `"A 2864-2026`")"
$ws.Range("Y18").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsynsmail/A 2864-2026 tillsynsbegäran mail.docx`", `"A 2864-2026`")"
$ws.Range("Z18").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/fåglar/A 2864-2026 prioriterade fågelarter.docx`", `"A 2864-2026`")"

# Row 19  <-  old row 20 data (A 61963-2025)
$ws.Range("A19").Value = "A 61963-2025"
$ws.Range("B19").Value = 46003.0
$ws.Range("G19").Value = 1.1
$ws.Range("H19").Value = 2.0
$ws.Range("I19").Value = 1.0
$ws.Range("J19").Value = 1.0
$ws.Range("K19").Value = 0.0
$ws.Range("L19").Value = 0.0
$ws.Range("M19").Value = 0.0
$ws.Range("N19").Value = 0.0
$ws.Range("O19").Value = 1.0
$ws.Range("P19").Value = 0.0
$ws.Range("Q19").Value = 3.0
$ws.Range("R19").Value = "Spillkråka`r`nMindre märgborre`r`nBlåsippa"
$ws.Range("S19").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/artfynd/A 61963-2025 artfynd.xlsx`", `"A 61963-2025`")"
$ws.Range("T19").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/kartor/A 61963-2025 karta.png`", `"A 61963-2025`")"
$ws.Range("V19").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomål/A 61963-2025 FSC-klagomål.docx`", `"A 61963-2025`")"
$ws.Range("W19").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomålsmail/A 61963-2025 FSC-klagomål mail.docx`", `"A 61963-2025`")"
$ws.Range("X19").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsyn/A 61963-2025 tillsynsbegäran.docx`", `"A 61963-2025`")"
$ws.Range("Y19").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsynsmail/A 61963-2025 tillsynsbegäran mail.docx`", `"A 61963-2025`")"
$ws.Range("Z19").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/fåglar/A 61963-2025 prioriterade fågelarter.docx`", `"A 61963-2025`")"

# Row 20  <-  old row 18 data (A 35242-2024)
$ws.Range("A20").Value = "A 35242-2024"
$ws.Range("B20").Value = 45530.55440972222
$ws.Range("G20").Value = 0.9
$ws.Range("H20").Value = 1.0
$ws.Range("I20").Value = 2.0
$ws.Range("J20").Value = 0.0
$ws.Range("K20").Value = 0.0
$ws.Range("L20").Value = 1.0
$ws.Range("M20").Value = 0.0
$ws.Range("N20").Value = 0.0
$ws.Range("O20").Value = 1.0
$ws.Range("P20").Value = 1.0
$ws.Range("Q20").Value = 3.0
$ws.Range("R20").Value = "Ryl`r`nGrönpyrola`r`nSkogsknipprot"
$ws.Range("S20").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/artfynd/A 35242-2024 artfynd.xlsx`", `"A 35242-2024`")"
$ws.Range("T20").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/kartor/A 35242-2024 karta.png`", `"A 35242-2024`")"
$ws.Range("V20").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomål/A 35242-2024 FSC-klagomål.docx`", `"A 35242-2024`")"
$ws.Range("W20").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomålsmail/A 35242-2024 FSC-klagomål mail.docx`", `"A 35242-2024`")"
$ws.Range("X20").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsyn/A 35242-2024 tillsynsbegäran.docx`", `"A 35242-2024`")"
$ws.Range("Y20").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsynsmail/A 35242-2024 tillsynsbegäran mail.docx`", `"A 35242-2024`")"
$ws.Range("Z20").ClearContents()

# Row 21  <-  old row 22 data (A 15600-2025)
$ws.Range("A21").Value = "A 15600-2025"
$ws.Range("B21").Value = 45747.0
$ws.Range("G21").Value = 1.1
$ws.Range("H21").Value = 0.0
$ws.Range("I21").Value = 2.0
$ws.Range("J21").Value = 0.0
$ws.Range("K21").Value = 0.0
$ws.Range("L21").Value = 0.0
$ws.Range("M21").Value = 0.0
$ws.Range("N21").Value = 0.0
$ws.Range("O21").Value = 0.0
$ws.Range("P21").Value = 0.0
$ws.Range("Q21").Value = 2.0
$ws.Range("R21").Value = "Murgröna`r`nScharlakansvårskål agg."
$ws.Range("S21").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/artfynd/A 15600-2025 artfynd.xlsx`", `"A 15600-2025`")"
$ws.Range("T21").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/kartor/A 15600-2025 karta.png`", `"A 15600-2025`")"
$ws.Range("V21").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomål/A 15600-2025 FSC-klagomål.docx`", `"A 15600-2025`")"
$ws.Range("W21").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomålsmail/A 15600-2025 FSC-klagomål mail.docx`", `"A 15600-2025`")"
$ws.Range("X21").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsyn/A 15600-2025 tillsynsbegäran.docx`", `"A 15600-2025`")"
$ws.Range("Y21").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsynsmail/A 15600-2025 tillsynsbegäran mail.docx`", `"A 15600-2025`")"

# Row 22  <-  old row 23 data (A 62231-2023)
$ws.Range("A22").Value = "A 62231-2023"
$ws.Range("B22").Value = 45267.0
$ws.Range("G22").Value = 1.6
$ws.Range("H22").Value = 1.0
$ws.Range("I22").Value = 1.0
$ws.Range("J22").Value = 0.0
$ws.Range("K22").Value = 0.0
$ws.Range("L22").Value = 0.0
$ws.Range("M22").Value = 0.0
$ws.Range("N22").Value = 0.0
$ws.Range("O22").Value = 0.0
$ws.Range("P22").Value = 0.0
$ws.Range("Q22").Value = 2.0
$ws.Range("R22").Value = "Murgröna`r`nBlåsippa"
$ws.Range("S22").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/artfynd/A 62231-2023 artfynd.xlsx`", `"A 62231-2023`")"
$ws.Range("T22").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/kartor/A 62231-2023 karta.png`", `"A 62231-2023`")"
$ws.Range("V22").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomål/A 62231-2023 FSC-klagomål.docx`", `"A 62231-2023`")"
$ws.Range("W22").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomålsmail/A 62231-2023 FSC-klagomål mail.docx`", `"A 62231-2023`")"
$ws.Range("X22").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsyn/A 62231-2023 tillsynsbegäran.docx`", `"A 62231-2023`")"
$ws.Range("Y22").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsynsmail/A 62231-2023 tillsynsbegäran mail.docx`", `"A 62231-2023`")"

# Row 23  <-  old row 21 data (A 40361-2023)
$ws.Range("A23").Value = "A 40361-2023"
$ws.Range("B23").Value = 45169.0
$ws.Range("G23").Value = 1.1
$ws.Range("H23").Value = 0.0
$ws.Range("I23").Value = 0.0
$ws.Range("J23").Value = 1.0
$ws.Range("K23").Value = 0.0
$ws.Range("L23").Value = 0.0
$ws.Range("M23").Value = 1.0
$ws.Range("N23").Value = 0.0
$ws.Range("O23").Value = 2.0
$ws.Range("P23").Value = 1.0
$ws.Range("Q23").Value = 2.0
$ws.Range("R23").Value = "Lundalm`r`nÄngsskära"
$ws.Range("S23").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/artfynd/A 40361-2023 artfynd.xlsx`", `"A 40361-2023`")"
$ws.Range("T23").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/kartor/A 40361-2023 karta.png`", `"A 40361-2023`")"
$ws.Range("V23").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomål/A 40361-2023 FSC-klagomål.docx`", `"A 40361-2023`")"
$ws.Range("W23").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomålsmail/A 40361-2023 FSC-klagomål mail.docx`", `"A 40361-2023`")"
$ws.Range("X23").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsyn/A 40361-2023 tillsynsbegäran.docx`", `"A 40361-2023`")"
$ws.Range("Y23").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsynsmail/A 40361-2023 tillsynsbegäran mail.docx`", `"A 40361-2023`")"

# Row 26  <-  old row 27 data (A 43702-2025)
$ws.Range("A26").Value = "A 43702-2025"
$ws.Range("B26").Value = 45912.40270833333
$ws.Range("G26").Value = 2.6
$ws.Range("H26").Value = 0.0
$ws.Range("I26").Value = 1.0
$ws.Range("J26").Value = 0.0
$ws.Range("K26").Value = 0.0
$ws.Range("L26").Value = 0.0
$ws.Range("M26").Value = 0.0
$ws.Range("N26").Value = 0.0
$ws.Range("O26").Value = 0.0
$ws.Range("P26").Value = 0.0
$ws.Range("Q26").Value = 1.0
$ws.Range("R26").Value = "Sårläka"
$ws.Range("S26").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/artfynd/A 43702-2025 artfynd.xlsx`", `"A 43702-2025`")"
$ws.Range("T26").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/kartor/A 43702-2025 karta.png`", `"A 43702-2025`")"
$ws.Range("V26").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomål/A 43702-2025 FSC-klagomål.docx`", `"A 43702-2025`")"
$ws.Range("W26").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomålsmail/A 43702-2025 FSC-klagomål mail.docx`", `"A 43702-2025`")"
$ws.Range("X26").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsyn/A 43702-2025 tillsynsbegäran.docx`", `"A 43702-2025`")"
$ws.Range("Y26").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsynsmail/A 43702-2025 tillsynsbegäran mail.docx`", `"A 43702-2025`")"

# Row 27  <-  old row 28 data (A 43704-2025)
$ws.Range("A27").Value = "A 43704-2025"
$ws.Range("B27").Value = 45912.40490740741
$ws.Range("G27").Value = 2.7
$ws.Range("H27").Value = 0.0
$ws.Range("I27").Value = 1.0
$ws.Range("J27").Value = 0.0
$ws.Range("K27").Value = 0.0
$ws.Range("L27").Value = 0.0
$ws.Range("M27").Value = 0.0
$ws.Range("N27").Value = 0.0
$ws.Range("O27").Value = 0.0
$ws.Range("P27").Value = 0.0
$ws.Range("Q27").Value = 1.0
$ws.Range("R27").Value = "Blomkålssvamp"
$ws.Range("S27").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/artfynd/A 43704-2025 artfynd.xlsx`", `"A 43704-2025`")"
$ws.Range("T27").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/kartor/A 43704-2025 karta.png`", `"A 43704-2025`")"
$ws.Range("V27").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomål/A 43704-2025 FSC-klagomål.docx`", `"A 43704-2025`")"
$ws.Range("W27").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomålsmail/A 43704-2025 FSC-klagomål mail.docx`", `"A 43704-2025`")"
$ws.Range("X27").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsyn/A 43704-2025 tillsynsbegäran.docx`", `"A 43704-2025`")"
$ws.Range("Y27").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsynsmail/A 43704-2025 tillsynsbegäran mail.docx`", `"A 43704-2025`")"

# Row 28  <-  old row 26 data (A 19459-2023)
$ws.Range("A28").Value = "A 19459-2023"
$ws.Range("B28").Value = 45049.0
$ws.Range("G28").Value = 1.6
$ws.Range("H28").Value = 0.0
$ws.Range("I28").Value = 0.0
$ws.Range("J28").Value = 0.0
$ws.Range("K28").Value = 1.0
$ws.Range("L28").Value = 0.0
$ws.Range("M28").Value = 0.0
$ws.Range("N28").Value = 0.0
$ws.Range("O28").Value = 1.0
$ws.Range("P28").Value = 1.0
$ws.Range("Q28").Value = 1.0
$ws.Range("R28").Value = "Luddvicker"
$ws.Range("S28").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/artfynd/A 19459-2023 artfynd.xlsx`", `"A 19459-2023`")"
$ws.Range("T28").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/kartor/A 19459-2023 karta.png`", `"A 19459-2023`")"
$ws.Range("V28").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomål/A 19459-2023 FSC-klagomål.docx`", `"A 19459-2023`")"
$ws.Range("W28").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/klagomålsmail/A 19459-2023 FSC-klagomål mail.docx`", `"A 19459-2023`")"
$ws.Range("X28").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsyn/A 19459-2023 tillsynsbegäran.docx`", `"A 19459-2023`")"
$ws.Range("Y28").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0840/tillsynsmail/A 19459-2023 tillsynsbegäran mail.docx`", `"A 19459-2023`")"
